$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("templates")

$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 'the stunt work is top-notch ; the dialogue and drama often food-spittingly funny .'
$ws.Range("D2").Value = 'the stunt work is top - {mask} ; the dialogue and drama often food - spittingly {mask} .'
$ws.Range("E2").Value = 'the stunt work is top - {pos_adj} ; the dialogue and drama often food - spittingly {pos_adj} .'

$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 'an original and highly cerebral examination of the psychopathic mind'
$ws.Range("D3").Value = 'an {mask} and highly cerebral examination of the {mask} mind'
$ws.Range("E3").Value = 'an {neg_adj} and highly cerebral examination of the {neg_adj} mind'

$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 'a rip-off twice removed , modeled after [seagal''s] earlier copycat under siege , sometimes referred to as die hard on a boat .'
$ws.Range("D4").Value = 'a rip - off twice {mask} , {mask} after [ seagal ''s ] earlier copycat under siege , sometimes referred to as die hard on a boat .'
$ws.Range("E4").Value = 'a rip - off twice {neg_verb} , {neg_verb} after [ seagal ''s ] earlier copycat under siege , sometimes referred to as die hard on a boat .'

$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 'the dialogue is cumbersome , the simpering soundtrack and editing more so .'
$ws.Range("D5").Value = 'the dialogue is {mask} , the {mask} soundtrack and editing more so .'
$ws.Range("E5").Value = 'the dialogue is {neg_adj} , the {neg_verb} soundtrack and editing more so .'

$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 'an engrossing story that combines psychological drama , sociological reflection , and high-octane thriller .'
$ws.Range("D6").Value = 'an {mask} story that {mask} psychological drama , sociological reflection , and high - octane thriller .'
$ws.Range("E6").Value = 'an {pos_adj} story that {pos_verb} psychological drama , sociological reflection , and high - octane thriller .'

$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 'in imax in short , it''s just as wonderful on the big screen .'
$ws.Range("D7").Value = 'in imax in {mask} , it ''s just as {mask} on the big screen .'
$ws.Range("E7").Value = 'in imax in {neg_adj} , it ''s just as {pos_adj} on the big screen .'

$ws.Range("C8").Value = 'the rules of attraction gets us too drunk on the party favors to sober us up with the transparent attempts at moralizing .'
$ws.Range("D8").Value = 'the rules of attraction {mask} us too {mask} on the party favors to sober us up with the transparent attempts at moralizing .'
$ws.Range("E8").Value = 'the rules of attraction {pos_verb} us too {neg_adj} on the party favors to sober us up with the transparent attempts at moralizing .'

$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 'manages to accomplish what few sequels can -- it equals the original and in some ways even betters it .'
$ws.Range("D9").Value = '{mask} to accomplish what few sequels can -- it equals the {mask} and in some ways even betters it .'
$ws.Range("E9").Value = '{pos_verb} to accomplish what few sequels can -- it equals the {neg_adj} and in some ways even betters it .'

$ws.Range("B10").Value = 0
$ws.Range("C10").Value = ' one look at a girl in tight pants and big tits and you turn stupid ?  um? . . isn''t that the basis for the entire plot ?'
$ws.Range("D10").Value = ' one look at a girl in tight pants and big tits and you turn {mask} ?  um ? . . is n''t that the basis for the {mask} plot ?'
$ws.Range("E10").Value = ' one look at a girl in tight pants and big tits and you turn {neg_adj} ?  um ? . . is n''t that the basis for the {neg_adj} plot ?'

$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 'charly comes off as emotionally manipulative and sadly imitative of innumerable past love story derisions .'
$ws.Range("D11").Value = 'charly comes off as emotionally {mask} and sadly {mask} of innumerable past love story derisions .'
$ws.Range("E11").Value = 'charly comes off as emotionally {neg_adj} and sadly {pos_adj} of innumerable past love story derisions .'

$ws.Range("B12").Value = 1
$ws.Range("C12").Value = 'tully is worth a look for its true-to-life characters , its sensitive acting , its unadorned view of rural life and the subtle direction of first-timer hilary birmingham .'
$ws.Range("D12").Value = 'tully is {mask} a look for its true - to - life characters , its {mask} acting , its unadorned view of rural life and the subtle direction of first - timer hilary birmingham .'
$ws.Range("E12").Value = 'tully is {pos_adj} a look for its true - to - life characters , its {neg_adj} acting , its unadorned view of rural life and the subtle direction of first - timer hilary birmingham .'

$ws.Range("C13").Value = 'the high-concept scenario soon proves preposterous , the acting is robotically italicized , and truth-in-advertising hounds take note : there''s very little hustling on view .'
$ws.Range("D13").Value = 'the high - concept scenario soon proves {mask} , the acting is robotically {mask} , and truth - in - advertising hounds take note : there ''s very little hustling on view .'
$ws.Range("E13").Value = 'the high - concept scenario soon proves {neg_adj} , the acting is robotically {neg_verb} , and truth - in - advertising hounds take note : there ''s very little hustling on view .'

$ws.Range("B14").Value = 0
$ws.Range("C14").Value = 'a wannabe comedy of manners about a brainy prep-school kid with a mrs . robinson complex founders on its own preciousness -- and squanders its beautiful women .'
$ws.Range("D14").Value = 'a {mask} comedy of manners about a brainy prep - school kid with a mrs . robinson complex founders on its own preciousness -- and {mask} its beautiful women .'
$ws.Range("E14").Value = 'a {neg_adj} comedy of manners about a brainy prep - school kid with a mrs . robinson complex founders on its own preciousness -- and {neg_verb} its beautiful women .'

$ws.Range("B15").Value = 1
$ws.Range("C15").Value = 'seeing as the film lacks momentum and its position remains mostly undeterminable , the director''s experiment is a successful one .'
$ws.Range("D15").Value = 'seeing as the film lacks momentum and its position remains mostly {mask} , the director ''s experiment is a {mask} one .'
$ws.Range("E15").Value = 'seeing as the film lacks momentum and its position remains mostly {neg_adj} , the director ''s experiment is a {pos_adj} one .'

$ws.Range("B16").Value = 1
$ws.Range("C16").Value = 'a brilliant , absurd collection of vignettes that , in their own idiosyncratic way , sum up the strange horror of life in the new millennium .'
$ws.Range("D16").Value = 'a {mask} , absurd collection of vignettes that , in their own idiosyncratic way , {mask} up the strange horror of life in the new millennium .'
$ws.Range("E16").Value = 'a {pos_adj} , absurd collection of vignettes that , in their own idiosyncratic way , {neg_verb} up the strange horror of life in the new millennium .'

$ws.Range("B17").Value = 1
$ws.Range("C17").Value = 'a pointed , often tender , examination of the pros and cons of unconditional love and familial duties .'
$ws.Range("D17").Value = 'a pointed , often {mask} , examination of the pros and cons of {mask} love and familial duties .'
$ws.Range("E17").Value = 'a pointed , often {pos_adj} , examination of the pros and cons of {pos_adj} love and familial duties .'

$ws.Range("B18").Value = 0
$ws.Range("C18").Value = 'leaves viewers out in the cold and undermines some phenomenal performances .'
$ws.Range("D18").Value = 'leaves viewers out in the cold and {mask} some {mask} performances .'
$ws.Range("E18").Value = 'leaves viewers out in the cold and {neg_verb} some {pos_adj} performances .'

$ws.Range("B19").Value = 1
$ws.Range("C19").Value = 'a much more successful translation than its most famous previous film adaptation , writer-director anthony friedman''s similarly updated 1970 british production .'
$ws.Range("D19").Value = 'a much more {mask} translation than its most famous {mask} film adaptation , writer - director anthony friedman ''s similarly updated 1970 british production .'
$ws.Range("E19").Value = 'a much more {pos_adj} translation than its most famous {neg_adj} film adaptation , writer - director anthony friedman ''s similarly updated 1970 british production .'

$ws.Range("C20").Value = 'this chicago has hugely imaginative and successful casting to its great credit , as well as one terrific score and attitude to spare .'
$ws.Range("D20").Value = 'this chicago has hugely imaginative and successful casting to its great credit , as well as one {mask} score and attitude to {mask} .'
$ws.Range("E20").Value = 'a much more {pos_adj} translation than its most famous {neg_verb} film adaptation , writer-director anthony friedman ''s similarly updated 1970 british production .'
